# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new column headers in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (bold, centered, bordered).
$ws.Range("C1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record is the same for every player on the roster (team record).
$wins = 82
$losses = 80
$ties = 0

for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins
    $ws.Cells.Item($row, 31).Value = $losses
    $ws.Cells.Item($row, 32).Value = $ties
}
